$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Glyphs")

$ws.Range("A50").Value = "g49"
$ws.Range("B50").Value = "k abbreviation"

$ws.Range("E4").Select()
